$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Query_1_DB" / "Query_2_DB" columns (E and F) used to hold a numeric
# database id ("1"); replace it with the database name ("Snowflake")
# for every data row (rows 2-18).
$ws.Range("E2:F18").Value = "Snowflake"

# Rows 12-16 previously used a slightly different border (no separate
# bottom rule) - normalize them to the same bordered style used by the
# rest of the column now that the content matches.
for ($r = 12; $r -le 16; $r++) {
    $rowRange = $ws.Range("E" + $r + ":F" + $r)
    $rowRange.Borders.Item(9).LineStyle = 1
    $rowRange.Borders.Item(9).Weight = 2
}

# Row 18 previously used a wrap-text style; normalize it too so it
# matches the plain bordered style used elsewhere.
$ws.Range("E18:F18").WrapText = $false

$ws.Range("D8").Select()
